$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: append " (transactions)" to the defect description; B3/C3 text unchanged
$ws.Range("A3").Value = "Backend and Frontend use different account transaction format (transactions)"
$ws.Range("B3").Value = "The Backend uses spaces, and the Frontend uses underscores"
$ws.Range("C3").Value = "Updated the Backend to read in input files with underscores"

# Row 6: updated defect text + new solution note
$ws.Range("A6").Value = "FrontEnd file output does match Backend file input 2"
$ws.Range("B6").Value = "Frontend doesn't have N or S on user accounts"
$ws.Range("C6").Value = "added to look like: 00001_John_Doe_____________A_00050.00_N"

# Row 7: new defect entry (accounts variant), reusing the same description as row 3/B3
$ws.Range("A7").Value = "Backend and Frontend use different account transaction format (accounts)"
$ws.Range("B7").Value = "The Backend uses spaces, and the Frontend uses underscores"

# Update the active selection to C7
$ws.Range("C7").Select()
